$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 was previously empty (sheet jumped from row 1 straight to row 3).
# Fill it in with a new hour-log entry, matching the date formatting used
# by the other entries in column A (numFmt "d-mmm", same as A3, A4, ...).
$ws.Range("A2").Value = 42207
$ws.Range("A2").NumberFormat = "d-mmm"
$ws.Range("B2").Value = 2

# Update the active cell / selection to A3, as in the target workbook.
$ws.Range("A3").Select()
